$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2; this pushes the existing rows 2..24 down to
# 3..25 (the "Fonte" hyperlink's cell content/style at O2 slides down to O3
# along with the rest of the row).
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the partial-year 2023 data (only Jan-Mar are
# available so far).
$ws.Range("A2").Value = 2023
$ws.Range("B2").Value = 1.1233
$ws.Range("C2").Value = 0.9181
$ws.Range("D2").Value = 1.1747000000000001

# The inserted row otherwise inherits styled-but-empty cells across E:N (from
# the column formatting); clear those so they don't materialize in the XML.
$ws.Range("E2:N2").Clear()

# Fill in the previously-missing December value for 2022 (now row 3).
$ws.Range("M3").Value = 1.1233

# The Hyperlinks collection still thinks the "Fonte" link is anchored at O2
# (its Range didn't automatically follow the row insert), so rebuild it
# pointing at the new location, O3, and restore the Hiperlink cell style.
$ws.Range("O2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("O3"), "https://www.valor.srv.br/indices/cdi.php")
$ws.Range("O3").Style = "Hiperlink"
